# "Add files via upload" — the uploaded Services.xlsx replaces the product
# image filename used by the two fiber-optic rows (E4 and E6, both driven by
# the same shared string) and leaves behind a new cursor / scroll position
# from whatever cell the author last clicked in Excel before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell whose "image" column held "fiber-optic.sedan.mp4" now holds
# "image_1353.jpg" (both Fiber Optic Lights Sedan rows share this string).
$ws.Range("E4").Value = "image_1353.jpg"
$ws.Range("E6").Value = "image_1353.jpg"

# Reflect the author's final on-screen selection: the view scrolled up a row
# (topLeftCell A5 -> A4) and the cursor ended on F4 (cell F4, after having
# E4 also highlighted).
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E4:F4").Select()
$ws.Range("F4").Activate()
